$wb = $excel.ActiveWorkbook

# Set the previous active sheet selection to B18 on sheet "1" before adding the new sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("B18").Select()

# Add new worksheet "11" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "11"

# Copy the header/scorer-column formatting (bold, centered, bordered) from sheet "1" column A
$ws1.Range("A1:A27").Copy()
$ws.Range("A1:A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate data column-by-column (A, then B, then C, then D, then E) to match
# the original shared-string insertion order used when this workbook was generated.
# Column A
$ws.Range("A1").Value = "Scorer"
$ws.Range("A2").Value = "абрахам"
$ws.Range("A3").Value = "аларио"
$ws.Range("A4").Value = "ассомбалонга"
$ws.Range("A5").Value = "балотелли"
$ws.Range("A6").Value = "бензема"
$ws.Range("A7").Value = "бла"
$ws.Range("A8").Value = "винд"
$ws.Range("A9").Value = "дембеле"
$ws.Range("A10").Value = "диаби"
$ws.Range("A11").Value = "дурсун"
$ws.Range("A12").Value = "дэвид"
$ws.Range("A13").Value = "зайц"
$ws.Range("A14").Value = "изидор"
$ws.Range("A15").Value = "иммобиле"
$ws.Range("A16").Value = "коло-муани"
$ws.Range("A17").Value = "коутиньо"
$ws.Range("A18").Value = "ляказетт"
$ws.Range("A19").Value = "милик"
$ws.Range("A20").Value = "обамеянг"
$ws.Range("A21").Value = "озил"
$ws.Range("A22").Value = "ракитич"
$ws.Range("A23").Value = "сака"
$ws.Range("A24").Value = "смит-роу"
$ws.Range("A25").Value = "трезеге"
$ws.Range("A26").Value = "ундер"
$ws.Range("A27").Value = "уткин"
$ws.Range("A28").Value = "чикаллеши"
$ws.Range("A29").Value = "эль-хадди"
$ws.Range("A30").Value = "торрес"

# Column B
$ws.Range("B1").Value = "Team"
$ws.Range("B2").Value = "рома"
$ws.Range("B3").Value = "байер"
$ws.Range("B4").Value = "адана"
$ws.Range("B5").Value = "адана демирспор"
$ws.Range("B6").Value = "реал мадрид"
$ws.Range("B7").Value = "нант"
$ws.Range("B8").Value = "вольфсбург"
$ws.Range("B9").Value = "барселона"
$ws.Range("B10").Value = "байер"
$ws.Range("B11").Value = "фенербахче"
$ws.Range("B12").Value = "лилль"
$ws.Range("B13").Value = "фенербахче"
$ws.Range("B14").Value = "локомотив"
$ws.Range("B15").Value = "лацио"
$ws.Range("B16").Value = "нант"
$ws.Range("B17").Value = "астон вилла"
$ws.Range("B18").Value = "арсенал"
$ws.Range("B19").Value = "марсель"
$ws.Range("B20").Value = "барселона"
$ws.Range("B21").Value = "фенербахче"
$ws.Range("B22").Value = "севилья"
$ws.Range("B23").Value = "арсенал"
$ws.Range("B24").Value = "арсенал"
$ws.Range("B25").Value = "башакшехир"
$ws.Range("B26").Value = "марсель"
$ws.Range("B27").Value = "ахмат"
$ws.Range("B28").Value = "коньяспор"
$ws.Range("B29").Value = "севилья"
$ws.Range("B30").Value = "барселона"

# Column C
$ws.Range("C1").Value = "Goals"
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 1
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 2
$ws.Range("C21").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 1
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("C28").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 1

# Column D
$ws.Range("D1").Value = "Picks"
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("D12").Value = 3
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("D21").Value = 3
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("D29").Value = 1
$ws.Range("D30").Value = 1

# Column E
$ws.Range("E1").Value = "Matchday"
$ws.Range("E2").Value = 11
$ws.Range("E3").Value = 11
$ws.Range("E4").Value = 11
$ws.Range("E5").Value = 11
$ws.Range("E6").Value = 11
$ws.Range("E7").Value = 11
$ws.Range("E8").Value = 11
$ws.Range("E9").Value = 11
$ws.Range("E10").Value = 11
$ws.Range("E11").Value = 11
$ws.Range("E12").Value = 11
$ws.Range("E13").Value = 11
$ws.Range("E14").Value = 11
$ws.Range("E15").Value = 11
$ws.Range("E16").Value = 11
$ws.Range("E17").Value = 11
$ws.Range("E18").Value = 11
$ws.Range("E19").Value = 11
$ws.Range("E20").Value = 11
$ws.Range("E21").Value = 11
$ws.Range("E22").Value = 11
$ws.Range("E23").Value = 11
$ws.Range("E24").Value = 11
$ws.Range("E25").Value = 11
$ws.Range("E26").Value = 11
$ws.Range("E27").Value = 11
$ws.Range("E28").Value = 11
$ws.Range("E29").Value = 11
$ws.Range("E30").Value = 11

$ws.Activate()
$ws.Range("H18").Select()

Write-Host "done"